$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("H2").Value = 2030
$ws.Range("D5").Value = 40

# Add new row 6
$ws.Range("C6").Value = "item"
$ws.Range("D6").Value = 11
$ws.Range("E6").Value = "Ether"
